$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The review originally in row 3 (nevilgreen@gmail.com / vikicrestina@gmail.com /
# 27/5/2019 11:25 / "עיצוב ידידותי מפתח ומלבב") is removed from its position;
# every row below shifts up by one, and that review is appended as the new
# last row (row 9). Also, the time for the review that becomes the new row 3
# (veredsnir12@gmail.com) is corrected from 28/5/2019 to 27/5/2019.

# 1) Remove row 3 - shifts rows 4..9 up to 3..8.
$ws.Rows.Item(3).Delete()

# 2) Re-create the removed review as the new last row (row 9), copying the
#    per-column formatting from cells that already carry the desired style.
$ws.Range("A2").Copy()
$ws.Range("A9").PasteSpecial(-4122)
$ws.Range("B2").Copy()
$ws.Range("B9").PasteSpecial(-4122)
$ws.Range("C2").Copy()
$ws.Range("C9").PasteSpecial(-4122)
$ws.Range("D3").Copy()
$ws.Range("D9").PasteSpecial(-4122)
$ws.Range("E2").Copy()
$ws.Range("E9").PasteSpecial(-4122)
$ws.Range("F2").Copy()
$ws.Range("F9").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

$ws.Range("A9").Value2 = "com.upstart42.dbcaptaincredit"
$ws.Range("B9").Value2 = "captain credit"
$ws.Range("C9").Value2 = "nevilgreen@gmail.com"
$ws.Range("D9").Value2 = "vikicrestina@gmail.com"
$ws.Range("E9").Value2 = "27/5/2019 11:25"
$ws.Range("F9").Value2 = "עיצוב ידידותי מפתח ומלבב"

# 3) Fix the time on the review that is now row 3.
$ws.Range("E3").Value2 = "27/5/2019 14:33"

# 4) Hyperlinks do not follow the row shift automatically, so rebuild the
#    whole hyperlink set to match the new row layout.
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("C2"), "mailto:leviadlevi22@gmail.com", [Type]::Missing, [Type]::Missing, "leviadlevi22@gmail.com") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D2"), "mailto:gazittalia1@gmail.com", [Type]::Missing, [Type]::Missing, "gazittalia1@gmail.com") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C3"), "mailto:veredsnir12@gmail.com", [Type]::Missing, [Type]::Missing, "veredsnir12@gmail.com") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D3"), "mailto:kevinkors122@gmail.com", [Type]::Missing, [Type]::Missing, "kevinkors122@gmail.com") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C4"), "mailto:stevewonder3001@gmail.com", [Type]::Missing, [Type]::Missing, "stevewonder3001@gmail.com") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D4"), "mailto:budoyoni@gmail.com", [Type]::Missing, [Type]::Missing, "budoyoni@gmail.com") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C5"), "mailto:stclerari834@gmail.com", [Type]::Missing, [Type]::Missing, "stclerari834@gmail.com") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C6"), "mailto:stcydouel274@gmail.com", [Type]::Missing, [Type]::Missing, "stcydouel274@gmail.com") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C7"), "mailto:kevinkors122@gmail.com", [Type]::Missing, [Type]::Missing, "kevinkors122@gmail.com") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D7"), "mailto:sinuspai@gmail.com", [Type]::Missing, [Type]::Missing, "sinuspai@gmail.com") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C8"), "mailto:halachme@gmail.com", [Type]::Missing, [Type]::Missing, "halachme@gmail.com") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C9"), "mailto:nevilgreen@gmail.com", [Type]::Missing, [Type]::Missing, "nevilgreen@gmail.com") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D9"), "mailto:vikicrestina@gmail.com", [Type]::Missing, [Type]::Missing, "vikicrestina@gmail.com") | Out-Null
